# LoginData.xlsx: add Sikuli login-data sheet (EmailAddress / Password) on
# Sheet2, mirroring the style already used on Sheet1, and make Sheet2 the
# active tab/sheet.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Selecting/activating Sheet2 is what flips workbook.xml's bookViews
# (activeTab) and drops tabSelected from Sheet1's sheetView.
$ws2.Activate()

$ws2.Range("A1").Value = "EmailAddress"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "swbtop@gmail.com"
$ws2.Range("B2").Value = "admin"

# Same cell style as Sheet1 (Text number format -> shared cellXfs entry).
$ws2.Range("A1:B2").NumberFormat = "@"

# Column A sized to fit its (now longer) contents, like Sheet1's columns.
$ws2.Columns("A:A").AutoFit()

# Portrait page orientation for the new sheet.
$ws2.PageSetup.Orientation = 1

# Leave the whole new table selected, anchored on the last entry.
$ws2.Range("A1:B2").Select()

Write-Output "Sheet2 populated with login data and activated"
